$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.381.44'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.066.35'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '234.53'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '57.07'
$ws.Range("E8").Value = '  -2.12%  '
$ws.Range("E9").Value = '  +2.47%  '
$ws.Range("D10").Value = '0.0777'
$ws.Range("E10").Value = '  +1.80%  '
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '2.371.32'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '14.33'
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '20.62'
$ws.Range("E14").Value = '  -3.41%  '
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '2.066.58'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("D18").Value = '37.314.69'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").Value = '6.23'
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '69.52'
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").Value = '226.09'
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  +1.41%  '
$ws.Range("E25").Value = '  -1.96%  '
$ws.Range("D26").Value = '167.97'
$ws.Range("E26").Value = '  +1.64%  '
$ws.Range("D27").Value = '8.82'
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("E29").Value = '  -7.12%  '
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").Value = '4.56'
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = '0.0617'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '4.54'
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").Value = '2.47'
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '5.61'
$ws.Range("E39").Value = '  -4.36%  '
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '1.493.39'
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").Value = '96.86'
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("D44").Value = '0.0212'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").Value = '4.18'
$ws.Range("E46").Value = '  -5.84%  '
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("E48").Value = '  -4.04%  '
$ws.Range("D49").Value = '7.22'
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("D50").Value = '2.97'
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").Value = '2.257.82'
$ws.Range("E51").Value = '  -0.04%  '
